$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new supplier data (replaces BTR/Hanun/Bantur)
$ws.Range("A2").Value = "SGS"
$ws.Range("B2").Value = "Ahmad"
$ws.Range("C2").Value = "Singosari"

# Remove row 3 (WGR/Pras/Wagir) entirely
$ws.Rows.Item(3).Delete()

# Update the selection to match the saved view state
$ws.Range("C10").Select()
